# Update LDLC prices history:
# Insert a new price-snapshot column right before the "nom" column (CH),
# shifting "nom" (was CH) to CI and "url_produit" (was CI) to CJ.
# The newly inserted column gets the new snapshot timestamp as its header
# (row 1) and, for every data row, the same price value that was most
# recently recorded (i.e. a copy of the former last price column, CG).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at CH; existing CH/CI (nom/url_produit) shift right
# to CI/CJ automatically, carrying over their values/styles.
$ws.Columns("CH").Insert()

# New header for the freshly inserted column CH (row 1).
$ws.Range("CH1").Value = "2026-01-31 14:14:23"

# Populate the new column's data rows (2-206) with the same values as the
# last existing price column (CG), which holds the latest known price for
# each product (keeps blanks blank where there is no price yet).
$ws.Range("CG2:CG206").Copy($ws.Range("CH2:CH206"))
